# Add the 32nd student ("Guzman, Maria Lilen") to the roster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new student's row (row 33 on the sheet, student #32)
$ws.Range("B33").Value = "Guzman, Maria Lilen"
$ws.Range("C33").Value = 1
$ws.Range("D33").Value = "lilenguzman2015@gmail.com"
$ws.Range("E33").Value = "lilenguzman01"

# Add a mailto hyperlink on the e-mail cell, matching the style used by the
# other rows in the "mail" column.
$ws.Hyperlinks.Add($ws.Range("D33"), "mailto:lilenguzman2015@gmail.com")

# Update the active cell/selection to reflect where the user ended up editing.
$ws.Range("E33").Select()
